$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowByAccount($acct) {
    $found = $ws.Columns.Item(1).Find($acct)
    if ($found) {
        return $found.Row()
    }
    return -1
}

function Set-AccountRow($row, $conta, $nome, $saldo) {
    $ws.Range("A$row").Value = "'" + $conta
    $ws.Range("B$row").Value = $nome
    $ws.Range("C$row").Value = $saldo
}

# 1. Insert a new row for MONICA right before the LEILA row (004208447).
$leilaRow = Get-RowByAccount("004208447")
$ws.Rows.Item($leilaRow).Insert()
Set-AccountRow $leilaRow "004387250" "MONICA" 500000

# 2. Replace the CLAUDIO row (008035153) with two rows: BERTILLA then ANTONIO.
$claudioRow = Get-RowByAccount("008035153")
$ws.Rows.Item($claudioRow).Delete()
$ws.Rows.Item($claudioRow).Insert()
$ws.Rows.Item($claudioRow).Insert()
Set-AccountRow $claudioRow "005064906" "BERTILLA" 40000
Set-AccountRow ($claudioRow + 1) "005000645" "ANTONIO" 24158.62

# 3. Insert a new row for JORGEANA right before the ANDRE row (005040864).
$andreRow = Get-RowByAccount("005040864")
$ws.Rows.Item($andreRow).Insert()
Set-AccountRow $andreRow "008002502" "JORGEANA" 11441.75

# 4. Remove the BLUEMETRIX row (001761119).
$bluemetrixRow = Get-RowByAccount("001761119")
$ws.Rows.Item($bluemetrixRow).Delete()
